$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Update "Marking" row correct-answer marks (B11): 3 -> 5
$ws.Range("B11").Value = 5

# Update "Total" row correct-answer marks (B12): 60 -> 100
$ws.Range("B12").Value = 100

# Update the total score fraction text (E12): "56/84" -> "100/140"
$ws.Range("E12").Value = "100/140"
